$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates - force text format so numeric-looking strings
# stay stored as text (matching the source data which is inlineStr).
$priceUpdates = @{
    "D2"  = "242.98"
    "D3"  = "23.10"
    "D4"  = "5.387"
    "D7"  = "6.479"
    "D8"  = "0.8093"
    "D9"  = "0.9084"
    "D10" = "0.1420"
    "D11" = "0.07415"
    "D12" = "0.03320"
    "D14" = "0.09336"
    "D15" = "3.847"
    "D16" = "0.001573"
    "D17" = "0.04632"
    "D19" = "0.006123"
    "D20" = "0.005031"
    "D21" = "0.0009842"
    "D22" = "0.00007796"
    "D24" = "3.612"
    "D27" = "0.1297"
    "D40" = "0.03889"
    "D41" = "0.006166"
    "D42" = "0.1070"
    "D44" = "0.007187"
    "D45" = "0.00005191"
    "D48" = "1.044"
    "D49" = "0.002260"
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $priceUpdates[$addr]
}

# Volume(1h) column (E) text updates
$ws.Range("E20").Value = "19HotbitTokenHTBBestin24h"
$ws.Range("E41").Value = "40KickTokenKICK"
